$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.516.22"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").Value = "3.877.18"
$ws.Range("E3").Value = "  -2.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.77"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.12"
$ws.Range("E6").Value = "  +5.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.670"
$ws.Range("E7").Value = "  -2.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").Value = "  +5.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.13"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.54"
$ws.Range("E13").Value = "  +5.00%  "

$ws.Range("D14").Value = "4.497.17"
$ws.Range("E14").Value = "  -2.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.18"
$ws.Range("E15").Value = "  +3.14%  "

$ws.Range("D16").Value = "3.872.49"
$ws.Range("E16").Value = "  -2.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.96"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("E18").Value = "  -3.61%  "

$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("D20").Value = "71.335.35"
$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.93"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.60"
$ws.Range("E23").Value = "  -2.06%  "

$ws.Range("E24").Value = "  -4.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.92"
$ws.Range("E25").Value = "  -3.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.69"
$ws.Range("E26").Value = "  +2.45%  "

$ws.Range("E27").Value = "  -5.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.57"
$ws.Range("E30").Value = "  +10.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.27"
$ws.Range("E31").Value = "  -3.20%  "

$ws.Range("E32").Value = "  -2.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.93"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("E35").Value = "  +10.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.80"
$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "630.67"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.29"
$ws.Range("E42").Value = "  -2.60%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.17"
$ws.Range("E43").Value = "  +19.18%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0474"
$ws.Range("E44").Value = "  -3.33%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("E45").Value = "  +7.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.24"
$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("E47").Value = "  -12.22%  "

$ws.Range("E48").Value = "  -3.57%  "

$ws.Range("D49").Value = "2.897.02"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000279"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.23"
$ws.Range("E51").Value = "  -4.86%  "
